# Performance.xlsx benchmark update
# Adds new "GCC 15.1" compiler-flag benchmark rows to the "Tex, Flags<0>"
# and "Flat, Flags<0>" sheets, then leaves the second sheet ("Flat, Flags<0>")
# as the active tab/selection the way the workbook was left after editing.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Tex, Flags<0>"
$ws2 = $wb.Worksheets.Item(2)   # "Flat, Flags<0>"

# --- Labels first, interleaved across both sheets so the shared-string
# table fills up in the same order the benchmark rows were typed in. ---
$ws1.Range("A27").Value = "GCC 15.1"
$ws2.Range("A20").Value = "GCC 15.1"

$ws2.Range("A21").Value = "Fast Fill"

$ws1.Range("A28").Value = "O3"
$ws2.Range("A22").Value = "O3"

$ws1.Range("A29").Value = "funswitch-loops"
$ws2.Range("A23").Value = "funswitch-loops"

$ws1.Range("A30").Value = "Hot Fn"
$ws2.Range("A24").Value = "Hot Fn"

$ws2.Range("A25").Value = "No align Fn"
$ws2.Range("A26").Value = "No sched inst"

# --- Sheet1 ("Tex, Flags<0>") data rows 27-30 ---
$ws1.Range("D27").Value = 110
$ws1.Range("E27").Formula = "=(D27/D`$2)-1"
$ws1.Range("F27").Formula = "=(D27/D26)-1"

$ws1.Range("D28").Value = 117
$ws1.Range("E28").Formula = "=(D28/D`$2)-1"
$ws1.Range("F28").Formula = "=(D28/D27)-1"
$ws1.Range("H28").Value = 17520

$ws1.Range("D29").Value = 110
$ws1.Range("E29").Formula = "=(D29/D`$2)-1"
$ws1.Range("F29").Formula = "=(D29/D28)-1"

$ws1.Range("D30").Value = 115
$ws1.Range("E30").Formula = "=(D30/D`$2)-1"
$ws1.Range("F30").Formula = "=(D30/D29)-1"
$ws1.Range("H30").Value = 14348

$ws1.Range("D27:D30").NumberFormat = "0"
$ws1.Range("E27:F30").NumberFormat = "0.0%"

# --- Sheet2 ("Flat, Flags<0>") data rows 20-26 ---
$ws2.Range("D20").Value = 444
$ws2.Range("E20").Formula = "=(D20/D`$2)-1"
$ws2.Range("F20").Formula = "=(D20/D19)-1"

$ws2.Range("D21").Value = 470
$ws2.Range("E21").Formula = "=(D21/D`$2)-1"
$ws2.Range("F21").Formula = "=(D21/D20)-1"

$ws2.Range("D22").Value = 689
$ws2.Range("E22").Formula = "=(D22/D`$2)-1"
$ws2.Range("F22").Formula = "=(D22/D21)-1"
$ws2.Range("H22").Value = 17520

$ws2.Range("D23").Value = 456
$ws2.Range("E23").Formula = "=(D23/D`$2)-1"
$ws2.Range("F23").Formula = "=(D23/D22)-1"
$ws2.Range("H23").Value = 14412

$ws2.Range("D24").Value = 499
$ws2.Range("E24").Formula = "=(D24/D`$2)-1"
$ws2.Range("F24").Formula = "=(D24/D23)-1"
$ws2.Range("H24").Value = 14652

$ws2.Range("D25").Value = 499
$ws2.Range("E25").Formula = "=(D25/D`$2)-1"
$ws2.Range("H25").Value = 14316

$ws2.Range("D26").Value = 499
$ws2.Range("E26").Formula = "=(D26/D`$2)-1"
$ws2.Range("H26").Value = 13952

$ws2.Range("E20:F24").NumberFormat = "0.0%"
$ws2.Range("E25:E26").NumberFormat = "0.0%"

# --- View state: sheet2 becomes the active tab, with its own selection;
# sheet1 keeps a selection too (it was the previously active sheet). ---
[void]$ws1.Range("A30").Select()
[void]$ws2.Activate()
[void]$ws2.Range("N13").Select()
